$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 445-447, pushing the existing rows (old 445-509) down to 448-512
$ws.Rows("445:447").Insert()

# Row 445
$ws.Cells.Item(445, 1).Value = 5
$ws.Cells.Item(445, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(445, 3).Value = "Maule"
$ws.Cells.Item(445, 4).Value = 45142
$ws.Cells.Item(445, 5).Value = 7
$ws.Cells.Item(445, 6).Value = "Fruta"
$ws.Cells.Item(445, 7).Value = 100102
$ws.Cells.Item(445, 8).Value = "Cítricos"
$ws.Cells.Item(445, 9).Value = 100102004
$ws.Cells.Item(445, 10).Value = "Mandarina"
$ws.Cells.Item(445, 11).Value = "Clemenuless"
$ws.Cells.Item(445, 12).Value = "Primera"
$ws.Cells.Item(445, 13).Value = 500
$ws.Cells.Item(445, 14).Value = 7000
$ws.Cells.Item(445, 15).Value = 7000
$ws.Cells.Item(445, 16).Value = 7000
$ws.Cells.Item(445, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(445, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(445, 19).Value = 700
$ws.Cells.Item(445, 20).Value = 10

# Row 446
$ws.Cells.Item(446, 1).Value = 5
$ws.Cells.Item(446, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(446, 3).Value = "Maule"
$ws.Cells.Item(446, 4).Value = 45142
$ws.Cells.Item(446, 5).Value = 7
$ws.Cells.Item(446, 6).Value = "Fruta"
$ws.Cells.Item(446, 7).Value = 100102
$ws.Cells.Item(446, 8).Value = "Cítricos"
$ws.Cells.Item(446, 9).Value = 100102004
$ws.Cells.Item(446, 10).Value = "Mandarina"
$ws.Cells.Item(446, 11).Value = "Murcott"
$ws.Cells.Item(446, 12).Value = "Primera"
$ws.Cells.Item(446, 13).Value = 280
$ws.Cells.Item(446, 14).Value = 10000
$ws.Cells.Item(446, 15).Value = 10000
$ws.Cells.Item(446, 16).Value = 10000
$ws.Cells.Item(446, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(446, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(446, 19).Value = 556
$ws.Cells.Item(446, 20).Value = 18

# Row 447
$ws.Cells.Item(447, 1).Value = 5
$ws.Cells.Item(447, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(447, 3).Value = "Maule"
$ws.Cells.Item(447, 4).Value = 45142
$ws.Cells.Item(447, 5).Value = 7
$ws.Cells.Item(447, 6).Value = "Fruta"
$ws.Cells.Item(447, 7).Value = 100102
$ws.Cells.Item(447, 8).Value = "Cítricos"
$ws.Cells.Item(447, 9).Value = 100102004
$ws.Cells.Item(447, 10).Value = "Mandarina"
$ws.Cells.Item(447, 11).Value = "Murcott"
$ws.Cells.Item(447, 12).Value = "Segunda"
$ws.Cells.Item(447, 13).Value = 200
$ws.Cells.Item(447, 14).Value = 8000
$ws.Cells.Item(447, 15).Value = 8000
$ws.Cells.Item(447, 16).Value = 8000
$ws.Cells.Item(447, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(447, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(447, 19).Value = 444
$ws.Cells.Item(447, 20).Value = 18
